# v 1.0.13 version prod
#
# The "Countries" sheet's Send_To / Send_To2 / Notify_To columns (J:L,
# rows 2-4) all share the same distribution-list string. Add the new
# recipient to that address, and drop the (no-op "No Fill") direct
# formatting that had been applied on top of those cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$distributionRange = $ws.Range("J2:L4")

# Append the new recipient to the shared distribution-list text.
$distributionRange.Value = "lucy.serrano@vcimentos.com;javier.martin@vcimentos.com;sat@rpatechnologies.es"

# Remove the explicit "No Fill" formatting that was sitting on these cells.
$distributionRange.Interior.Pattern = -4142  # xlPatternNone

# Leave the selection on J2 (matching the saved view state).
$ws.Activate()
$ws.Range("J2").Select()
